# Update error norm values for 'first_half_of_grid' and 'whole_grid' sheets
# and remove the 100000-column (F) which is no longer computed (Godunov method not working).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("first_half_of_grid")

$ws.Range("C2").Value2 = 0.3164403483074749
$ws.Range("D2").Value2 = 0.3233005541585267
$ws.Range("E2").Value2 = 0.3239753483761879
$ws.Range("C4").Value2 = 0.272518181771069
$ws.Range("D4").Value2 = 0.2810428891827828
$ws.Range("E4").Value2 = 0.2818495955627683
$ws.Range("C6").Value2 = 0.8346892101176664
$ws.Range("D6").Value2 = 0.8293338358937191
$ws.Range("E6").Value2 = 0.8287942655498161
$ws.Range("C8").Value2 = 2.237571161292486
$ws.Range("D8").Value2 = 7.229220162618181
$ws.Range("E8").Value2 = 22.90851657740766
$ws.Range("C10").Value2 = 1.926994543269511
$ws.Range("D10").Value2 = 6.284310048056426
$ws.Range("E10").Value2 = 19.92977602971193
$ws.Range("C12").Value2 = 5.902144006574449
$ws.Range("D12").Value2 = 18.54446833099011
$ws.Range("E12").Value2 = 58.60460453787991
$ws.Range("C14").Value2 = 0.09678932251842197
$ws.Range("D14").Value2 = 0.09991836081082778
$ws.Range("E14").Value2 = 0.1002316967432362
$ws.Range("C16").Value2 = 0.08113919979398444
$ws.Range("D16").Value2 = 0.08442616689516634
$ws.Range("E16").Value2 = 0.08474604094514103
$ws.Range("C18").Value2 = 0.3883858817270134
$ws.Range("D18").Value2 = 0.3842082113089135
$ws.Range("E18").Value2 = 0.3837912283306051
$ws.Range("C20").Value2 = 0.6968698230344685
$ws.Range("D20").Value2 = 0.6968698230344741
$ws.Range("E20").Value2 = 0.6968698230344821
$ws.Range("C22").Value2 = 0.6831492406486595
$ws.Range("D22").Value2 = 0.706320681331343
$ws.Range("E22").Value2 = 0.70874530287126
$ws.Range("C24").Value2 = 1.000000000000001
$ws.Range("D24").Value2 = 1.000000000000008
$ws.Range("E24").Value2 = 1.000000000000026

# Remove column F entirely (header + data + "decreasing rate" formulas)
$ws.Range("F1:F25").ClearContents()

$ws = $wb.Worksheets.Item("whole_grid")

$ws.Range("C2").Value2 = 0.5816297491623608
$ws.Range("D2").Value2 = 0.5860092954036585
$ws.Range("E2").Value2 = 0.5863339232404196
$ws.Range("C4").Value2 = 0.8973666580279368
$ws.Range("D4").Value2 = 0.9004080863272017
$ws.Range("E4").Value2 = 0.9006386914851281
$ws.Range("C6").Value2 = 0.700176950957232
$ws.Range("D6").Value2 = 0.7033808793887799
$ws.Range("E6").Value2 = 0.7034164424818655
$ws.Range("C8").Value2 = 5.816297491623607
$ws.Range("D8").Value2 = 18.53124103506002
$ws.Range("E8").Value2 = 58.63339232404195
$ws.Range("C10").Value2 = 8.973666580279367
$ws.Range("D10").Value2 = 28.47340376427471
$ws.Range("E10").Value2 = 90.06386914851281
$ws.Range("C12").Value2 = 7.00176950957232
$ws.Range("D12").Value2 = 22.24285641480728
$ws.Range("E12").Value2 = 70.34164424818655
$ws.Range("C14").Value2 = 0.4736624588108304
$ws.Range("D14").Value2 = 0.4786196686959243
$ws.Range("E14").Value2 = 0.479034569716764
$ws.Range("C16").Value2 = 0.6825826502368658
$ws.Range("D16").Value2 = 0.6873588477355467
$ws.Range("E16").Value2 = 0.6877835707369511
$ws.Range("C18").Value2 = 0.5545029375909902
$ws.Range("D18").Value2 = 0.5586723407762884
$ws.Range("E18").Value2 = 0.5587190841093198
$ws.Range("C20").Value2 = 0.9000000000000011
$ws.Range("D20").Value2 = 0.900000000000008
$ws.Range("E20").Value2 = 0.9000000000000198
$ws.Range("C22").Value2 = 1.699466079288994
$ws.Range("D22").Value2 = 1.73363170599307
$ws.Range("E22").Value2 = 1.734426288713234
$ws.Range("C24").Value2 = 1.000000000000001
$ws.Range("D24").Value2 = 1.000000000000009
$ws.Range("E24").Value2 = 1.000000000000042

# Remove column F entirely (header + data + "decreasing rate" formulas)
$ws.Range("F1:F25").ClearContents()

